$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 and A2 share the same style: bold font, thin box border, center/top alignment
$styledRange = $ws.Range("B1,A2")
$styledRange.Font.Bold = $true
$styledRange.Borders.LineStyle = 1
$styledRange.Borders.Weight = 2
$styledRange.HorizontalAlignment = -4108
$styledRange.VerticalAlignment = -4160

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0

# B2: shared string "disconnected_elements", default style
$ws.Range("B2").Value = "disconnected_elements"
